$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Sheet1" to "cuentas"
$ws.Name = "cuentas"

# Reset the active selection back to A1 (was A2)
$ws.Range("A1").Select()
